$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2..8) {
    $ws.Cells.Item($row, 3).Value = 45204
}
